# Commit: "Stated making OOP version of DD20parser"
# Rename the two mapping-sample values (dash -> underscore before trailing "1")
# and move the active selection to M3:M4, as captured in the last session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "E_EEE-FFF_1"
$ws.Range("B2").Value = "E_EEEV-FFF_1"

$null = $ws.Range("M3:M4").Select()
